$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '97.011.71'
$ws.Range("E2").Value = '  +1.97%  '
$ws.Range("D3").Value = '3.590.15'
$ws.Range("E3").Value = '  +0.45%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '240.42'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.49%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '659.90'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.44%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.69'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +14.84%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.426'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +6.77%  '
$ws.Range("E9").Value = '  -0.02%  '
$ws.Range("E10").Value = '  +3.97%  '
$ws.Range("D11").Value = '3.588.62'
$ws.Range("E11").Value = '  +0.49%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '44.12'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.87%  '
$ws.Range("E13").Value = '  +0.50%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.41'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.53%  '
$ws.Range("D15").Value = '4.259.65'
$ws.Range("E15").Value = '  -0.49%  '
$ws.Range("D16").Value = '96.881.69'
$ws.Range("E16").Value = '  +1.94%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000260'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.66%  '
$ws.Range("B18").Value = 'Polkadot'
$ws.Range("C18").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '8.60'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +9.06%  '
$ws.Range("B19").Value = 'WrappedEther'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D19").Value = '3.585.00'
$ws.Range("E19").Value = '  +0.48%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.72'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.31%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '18.03'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.43%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.523'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +9.43%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.52'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.44%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '514.36'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.37%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000203'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.45%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.89'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.40%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '101.03'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +6.00%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '13.04'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.36%  '
$ws.Range("D29").Value = '3.783.28'
$ws.Range("E29").Value = '  +0.59%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.159'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +11.41%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.01'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.22%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '11.85'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.36%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.997'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.47%  '
$ws.Range("E34").Value = '  +3.71%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.01'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.47%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '31.67'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.35%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '624.94'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +6.92%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.88'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +4.27%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.566'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.08%  '
$ws.Range("E40").Value = '  +1.29%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.96'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +10.03%  '
$ws.Range("E42").Value = '  +2.33%  '
$ws.Range("E43").Value = '  -0.05%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.923'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.93%  '
$ws.Range("E45").Value = '  +5.34%  '
$ws.Range("E46").Value = '  +6.92%  '
$ws.Range("E47").Value = '  -0.02%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '23.60'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.90%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.409'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +32.80%  '
$ws.Range("E50").Value = '  +4.82%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '33.14'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.70%  '
